$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Insert a new row before row 9 (current estimate_params row), shifting rows 9-17 down to 10-18
$ws.Rows.Item(9).Insert()

# Populate new row 9 with the L_curve parameter
$ws.Range("A9").Value2 = "L_curve"
$ws.Range("B9").Value2 = 0

# Rename the "Model" label (row 8) to "production_function"
$ws.Range("A8").Value2 = "production_function"

# Update the selection on this sheet and make it the active sheet/tab
$ws.Activate()
$ws.Range("B12:C13").Select()

Write-Host "done"
